$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need NumberFormat forced to
# text ("@") before assignment so Excel stores them as text (matching the
# source workbook, which keeps every Price/Volume cell as text), then the
# format is restored to General so no stray formatting diff is introduced.
$textSafeCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D22", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D47", "D48", "D50", "D51")
foreach ($cellRef in $textSafeCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.473.32'
$ws.Range("E2").Value = '  -2.79%  '
$ws.Range("D3").Value = '1.775.41'
$ws.Range("E3").Value = '  -1.75%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '305.06'
$ws.Range("E6").Value = '  -1.77%  '
$ws.Range("D7").Value = '0.4289'
$ws.Range("E7").Value = '  +1.87%  '
$ws.Range("E8").Value = '  +2.38%  '
$ws.Range("D9").Value = '0.07161'
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("D10").Value = '0.8473'
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("D11").Value = '20.58'
$ws.Range("E11").Value = '  +2.18%  '
$ws.Range("D12").Value = '1.774.69'
$ws.Range("E12").Value = '  -3.02%  '
$ws.Range("D13").Value = '5.260'
$ws.Range("E13").Value = '  -1.24%  '
$ws.Range("D14").Value = '6.451'
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").Value = '0.06873'
$ws.Range("E15").Value = '  -0.65%  '
$ws.Range("D17").Value = '79.01'
$ws.Range("E17").Value = '  -2.64%  '
$ws.Range("D18").Value = '0.000008710'
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.35%  '
$ws.Range("D20").Value = '14.99'
$ws.Range("E20").Value = '  -0.69%  '
$ws.Range("D21").Value = '26.490.57'
$ws.Range("E21").Value = '  -3.10%  '
$ws.Range("D22").Value = '5.125'
$ws.Range("E22").Value = '  +1.02%  '
$ws.Range("E23").Value = '  +1.97%  '
$ws.Range("D24").Value = '2.009.44'
$ws.Range("E24").Value = '  -2.28%  '
$ws.Range("D25").Value = '152.20'
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("E26").Value = '  -4.61%  '
$ws.Range("D27").Value = '18.05'
$ws.Range("E27").Value = '  -0.85%  '
$ws.Range("D28").Value = '5.083'
$ws.Range("E28").Value = '  +0.53%  '
$ws.Range("D29").Value = '113.91'
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("D30").Value = '1.796'
$ws.Range("E30").Value = '  +4.67%  '
$ws.Range("D31").Value = '0.08898'
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").Value = '0.7284'
$ws.Range("E32").Value = '  -1.57%  '
$ws.Range("D33").Value = '1.124'
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("D34").Value = '4.330'
$ws.Range("E34").Value = '  -3.00%  '
$ws.Range("B35").Value = 'Frax'
$ws.Range("C35").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D35").Value = '1.001'
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.741'
$ws.Range("E36").Value = '  -7.07%  '
$ws.Range("D37").Value = '1.098'
$ws.Range("E37").Value = '  +2.74%  '
$ws.Range("D38").Value = '0.05157'
$ws.Range("E38").Value = '  -1.07%  '
$ws.Range("D39").Value = '0.01886'
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("D40").Value = '0.4933'
$ws.Range("E40").Value = '  -0.98%  '
$ws.Range("E41").Value = '  -1.55%  '
$ws.Range("D42").Value = '2.652'
$ws.Range("E42").Value = '  -3.92%  '
$ws.Range("D43").Value = '6.325'
$ws.Range("E43").Value = '  +0.68%  '
$ws.Range("D44").Value = '8.025'
$ws.Range("E44").Value = '  -2.17%  '
$ws.Range("D45").Value = '105.22'
$ws.Range("E45").Value = '  +0.37%  '
$ws.Range("E46").Value = '  -1.09%  '
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").Value = '1.639'
$ws.Range("E48").Value = '  +2.77%  '
$ws.Range("E49").Value = '  -3.19%  '
$ws.Range("D50").Value = '0.4494'
$ws.Range("E50").Value = '  -1.75%  '
$ws.Range("D51").Value = '1.721'
$ws.Range("E51").Value = '  +3.18%  '

foreach ($cellRef in $textSafeCells) {
    $ws.Range($cellRef).NumberFormat = "General"
}
